# "add cell var callback"
# Adds a new row (row 10) below the existing grade-table block containing
# three template placeholder cells used by a cell-level callback:
#   A10 = "[A10]"  (italic)                     - the cell address marker
#   B10 = "x"       (yellow highlight fill)      - the callback return value
#   C10 = "{x}"     (default formatting)         - the template placeholder

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new row 10 content -----------------------------------------------
$ws.Range("A10").Value = "[A10]"
$ws.Range("B10").Value = "x"
$ws.Range("C10").Value = "{x}"

# A10: italicize, reusing the workbook's existing CJK-aware body font
# (copy format from a cell governed by row 2's custom font so the new
# italic font entry matches the sheet's existing font family/scheme).
$ws.Range("D2").Copy() | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").Font.Italic = $true
$excel.CutCopyMode = $false

# B10: highlight with a solid yellow fill
$ws.Range("B10").Interior.Color = 65535

# move the active selection to the newly added cell, like a user would
# after typing the last entry in the row
[void]$ws.Range("C10").Select()
